$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("table_definitions")
$ws1.Activate()
$ws1.Range("F4").Select()
Write-Host "ok"
